# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect the refreshed scrape output (gh-pages data at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$sheetExhibition.Range("F2").Value  = 828
$sheetExhibition.Range("F11").Value = 572
$sheetExhibition.Range("F13").Value = 13538
$sheetExhibition.Range("F15").Value = 20
$sheetExhibition.Range("F17").Value = 5566
$sheetExhibition.Range("F18").Value = 5589
$sheetExhibition.Range("F19").Value = 60

# 全部类型 sheet updates (mirrors the same events)
$sheetAllTypes.Range("F2").Value  = 828
$sheetAllTypes.Range("F33").Value = 572
$sheetAllTypes.Range("F35").Value = 13538
$sheetAllTypes.Range("F37").Value = 20
$sheetAllTypes.Range("F40").Value = 5566
$sheetAllTypes.Range("F41").Value = 5589
$sheetAllTypes.Range("F42").Value = 60
